$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.299.61"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.943.98"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'569.69"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "'159.37"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "2.939.10"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "'34.51"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "65.294.92"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "3.431.13"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "'7.03"
$ws.Range("D19").Value = "2.938.31"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("E20").Value = "  +10.98%  "
$ws.Range("D21").Value = "'445.22"
$ws.Range("E21").Value = "  -4.47%  "
$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "'82.46"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'12.16"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'8.03"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").Value = "'27.32"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "'5.76"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'0.971"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").Value = "'49.55"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "'44.31"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -8.91%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E42").Value = "  -7.52%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "'8.52"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "'383.48"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "2.697.34"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").Value = "'134.03"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E50").Value = "  +4.87%  "
$ws.Range("D51").Value = "'23.49"
$ws.Range("E51").Value = "  -0.63%  "
